$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value2 = 1513.5385
$ws.Range("I4").Value2 = 1445
$ws.Range("K4").Value2 = 1445
$ws.Range("M4").Value2 = -1331
$ws.Range("H12").Value2 = 350
$ws.Range("J12").Value2 = 0
$ws.Range("L12").Value2 = 0
$ws.Range("N12").ClearContents()
$ws.Range("H19").Value2 = 7144781
$ws.Range("J19").Value2 = 10002018
$ws.Range("L19").Value2 = 10002018
$ws.Range("N19").Value2 = -10002368
$ws.Range("H64").Value2 = 4026
$ws.Range("J64").Value2 = 4052.125
$ws.Range("L64").Value2 = 4052.125
$ws.Range("N64").Value2 = -4548.125
$ws.Range("H67").Value2 = 4026
$ws.Range("J67").Value2 = 4052.125
$ws.Range("L67").Value2 = 4052.125
$ws.Range("N67").Value2 = -5768.125
$ws.Range("H74").Value2 = 7678.4
$ws.Range("I74").Value2 = 7598.25
$ws.Range("K74").Value2 = 7598.25
$ws.Range("M74").Value2 = -6662.25
$ws.Range("H77").Value2 = 7678.4
$ws.Range("I77").Value2 = 7598.25
$ws.Range("K77").Value2 = 37991.25
$ws.Range("M77").Value2 = -33311.25
$ws.Range("H98").Value2 = 1985.6666
$ws.Range("I98").Value2 = 1842.2941
$ws.Range("K98").Value2 = 1842.2941
$ws.Range("M98").Value2 = -344.2941000000001
$ws.Range("H113").Value2 = 6495
$ws.Range("I113").Value2 = 6495
$ws.Range("J113").Value2 = 0
$ws.Range("K113").Value2 = 6495
$ws.Range("L113").Value2 = 0
$ws.Range("M113").Value2 = -3241
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value2 = 1985.6666
$ws.Range("I122").Value2 = 1842.2941
$ws.Range("K122").Value2 = 5526.8823
$ws.Range("M122").Value2 = -3076.8823
$ws.Range("H132").Value2 = 4537.256
$ws.Range("I132").Value2 = 2734.7
$ws.Range("J132").Value2 = 8697
$ws.Range("K132").Value2 = 8204.099999999999
$ws.Range("L132").Value2 = 26091
$ws.Range("M132").Value2 = -5674.099999999999
$ws.Range("N132").Value2 = -31151
$ws.Range("H135").Value2 = 3725.75
$ws.Range("J135").Value2 = 7776.8
$ws.Range("L135").Value2 = 69991.2
$ws.Range("N135").Value2 = -75061.2
$ws.Range("H137").Value2 = 2008.421
$ws.Range("I137").Value2 = 1963.625
$ws.Range("K137").Value2 = 5890.875
$ws.Range("M137").Value2 = -3340.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 3215.1516
$ws.Range("I32").Value2 = 2362.1206
$ws.Range("K32").Value2 = 2362.1206
$ws.Range("M32").Value2 = -2075.1206
$ws.Range("H61").Value2 = 20910742
$ws.Range("I61").Value2 = 26251772
$ws.Range("J61").Value2 = 6667999
$ws.Range("K61").Value2 = 26251772
$ws.Range("L61").Value2 = 6667999
$ws.Range("M61").Value2 = -26251560
$ws.Range("N61").Value2 = -6668423
$ws.Range("H74").Value2 = 2785.5789
$ws.Range("I74").Value2 = 1582.25
$ws.Range("J74").Value2 = 4848.4287
$ws.Range("K74").Value2 = 1582.25
$ws.Range("L74").Value2 = 4848.4287
$ws.Range("M74").Value2 = -708.25
$ws.Range("N74").Value2 = -6596.4287
$ws.Range("H77").Value2 = 2785.5789
$ws.Range("I77").Value2 = 1582.25
$ws.Range("J77").Value2 = 4848.4287
$ws.Range("K77").Value2 = 7911.25
$ws.Range("L77").Value2 = 24242.1435
$ws.Range("M77").Value2 = -3543.25
$ws.Range("N77").Value2 = -32978.14350000001
$ws.Range("H110").Value2 = 2948.1667
$ws.Range("I110").Value2 = 897
$ws.Range("J110").Value2 = 4999.3335
$ws.Range("K110").Value2 = 897
$ws.Range("L110").Value2 = 4999.3335
$ws.Range("M110").Value2 = 1148
$ws.Range("N110").Value2 = -9089.333500000001
$ws.Range("H132").Value2 = 2385994.2
$ws.Range("I132").Value2 = 4587.6177
$ws.Range("K132").Value2 = 13762.8531
$ws.Range("M132").Value2 = -11232.8531
$ws.Range("H136").Value2 = 20910742
$ws.Range("I136").Value2 = 26251772
$ws.Range("J136").Value2 = 6667999
$ws.Range("K136").Value2 = 78755316
$ws.Range("L136").Value2 = 20003997
$ws.Range("M136").Value2 = -78752766
$ws.Range("N136").Value2 = -20009097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 52695.715
$ws.Range("I86").Value2 = 97232.71000000001
$ws.Range("K86").Value2 = 97232.71000000001
$ws.Range("M86").Value2 = -96109.71000000001
$ws.Range("H89").Value2 = 52695.715
$ws.Range("I89").Value2 = 97232.71000000001
$ws.Range("K89").Value2 = 486163.55
$ws.Range("M89").Value2 = -480547.55
$ws.Range("H134").Value2 = 5265871.5
$ws.Range("I134").Value2 = 2648.7144
$ws.Range("K134").Value2 = 7946.1432
$ws.Range("M134").Value2 = -5411.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 5897155.5
$ws.Range("J16").Value2 = 2999
$ws.Range("L16").Value2 = 2999
$ws.Range("N16").Value2 = -3573
$ws.Range("H31").Value2 = 38465108
$ws.Range("I31").Value2 = 66670164
$ws.Range("J31").Value2 = 3670.3635
$ws.Range("K31").Value2 = 66670164
$ws.Range("L31").Value2 = 3670.3635
$ws.Range("M31").Value2 = -66669869
$ws.Range("N31").Value2 = -4260.363499999999
$ws.Range("H34").Value2 = 38465108
$ws.Range("I34").Value2 = 66670164
$ws.Range("J34").Value2 = 3670.3635
$ws.Range("K34").Value2 = 66670164
$ws.Range("L34").Value2 = 3670.3635
$ws.Range("M34").Value2 = -66669962
$ws.Range("N34").Value2 = -4074.3635
$ws.Range("H107").Value2 = 1175.1333
$ws.Range("I107").Value2 = 892.61536
$ws.Range("J107").Value2 = 3011.5
$ws.Range("K107").Value2 = 892.61536
$ws.Range("L107").Value2 = 3011.5
$ws.Range("M107").Value2 = 1027.38464
$ws.Range("N107").Value2 = -6851.5
$ws.Range("H113").Value2 = 5897155.5
$ws.Range("J113").Value2 = 2999
$ws.Range("L113").Value2 = 2999
$ws.Range("N113").Value2 = -7339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value2 = 5095.5
$ws.Range("I11").Value2 = 1818
$ws.Range("J11").Value2 = 9684
$ws.Range("K11").Value2 = 5454
$ws.Range("L11").Value2 = 29052
$ws.Range("M11").Value2 = -5314
$ws.Range("N11").Value2 = -29332
$ws.Range("H107").Value2 = 10112573
$ws.Range("J107").Value2 = 13001797
$ws.Range("L107").Value2 = 39005391
$ws.Range("N107").Value2 = -39009231
$ws.Range("H132").Value2 = 1827
$ws.Range("I132").Value2 = 1859.5
$ws.Range("J132").Value2 = 1794.5
$ws.Range("K132").Value2 = 16735.5
$ws.Range("L132").Value2 = 16150.5
$ws.Range("M132").Value2 = -14205.5
$ws.Range("N132").Value2 = -21210.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 6618.324
$ws.Range("I70").Value2 = 4378.6
$ws.Range("K70").Value2 = 4378.6
$ws.Range("M70").Value2 = -4108.6
$ws.Range("H73").Value2 = 6618.324
$ws.Range("I73").Value2 = 4378.6
$ws.Range("K73").Value2 = 4378.6
$ws.Range("M73").Value2 = -3442.6
$ws.Range("H98").Value2 = 11471
$ws.Range("J98").Value2 = 11471
$ws.Range("L98").Value2 = 11471
$ws.Range("N98").Value2 = -17461
$ws.Range("H102").Value2 = 2669.5625
$ws.Range("I102").Value2 = 2626.3215
$ws.Range("K102").Value2 = 2626.3215
$ws.Range("M102").Value2 = -1004.3215
$ws.Range("H107").Value2 = 1297.3914
$ws.Range("I107").Value2 = 1298.2106
$ws.Range("J107").Value2 = 1293.5
$ws.Range("K107").Value2 = 1298.2106
$ws.Range("L107").Value2 = 1293.5
$ws.Range("M107").Value2 = 621.7893999999999
$ws.Range("N107").Value2 = -5133.5
$ws.Range("H113").Value2 = 1686527.1
$ws.Range("I113").Value2 = 3359.4
$ws.Range("J113").Value2 = 3089167
$ws.Range("K113").Value2 = 3359.4
$ws.Range("L113").Value2 = 3089167
$ws.Range("M113").Value2 = -1189.4
$ws.Range("N113").Value2 = -3093507
$ws.Range("H126").Value2 = 6073
$ws.Range("I126").Value2 = 5996.4
$ws.Range("J126").Value2 = 6168.75
$ws.Range("K126").Value2 = 17989.2
$ws.Range("L126").Value2 = 18506.25
$ws.Range("M126").Value2 = -15519.2
$ws.Range("N126").Value2 = -23446.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 1855.75
$ws.Range("I46").Value2 = 1499.3334
$ws.Range("J46").Value2 = 2925
$ws.Range("K46").Value2 = 1499.3334
$ws.Range("L46").Value2 = 2925
$ws.Range("M46").Value2 = -1311.3334
$ws.Range("N46").Value2 = -3301
$ws.Range("H61").Value2 = 7599
$ws.Range("I61").Value2 = 1679.6
$ws.Range("K61").Value2 = 1679.6
$ws.Range("M61").Value2 = -1477.6
$ws.Range("H113").Value2 = 7599
$ws.Range("I113").Value2 = 1679.6
$ws.Range("K113").Value2 = 1679.6
$ws.Range("M113").Value2 = 490.4000000000001
$ws.Range("H122").Value2 = 3654.2
$ws.Range("I122").Value2 = 3336.6287
$ws.Range("J122").Value2 = 5877.2
$ws.Range("K122").Value2 = 10009.8861
$ws.Range("L122").Value2 = 17631.6
$ws.Range("M122").Value2 = -7559.8861
$ws.Range("N122").Value2 = -22531.6
$ws.Range("H132").Value2 = 4003
$ws.Range("I132").Value2 = 3060.389
$ws.Range("K132").Value2 = 9181.167000000001
$ws.Range("M132").Value2 = -6651.167000000001
$ws.Range("H136").Value2 = 5727.4287
$ws.Range("I136").Value2 = 5727.4287
$ws.Range("K136").Value2 = 17182.2861
$ws.Range("M136").Value2 = -14632.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value2 = 99992
$ws.Range("J16").Value2 = 99992
$ws.Range("L16").Value2 = 99992
$ws.Range("N16").Value2 = -100576
$ws.Range("H113").Value2 = 671.8182
$ws.Range("I113").Value2 = 568.25
$ws.Range("J113").Value2 = 796.1
$ws.Range("K113").Value2 = 1704.75
$ws.Range("L113").Value2 = 2388.3
$ws.Range("M113").Value2 = 465.25
$ws.Range("N113").Value2 = -6728.3
$ws.Range("H132").Value2 = 284920.38
$ws.Range("I132").Value2 = 6907.6772
$ws.Range("K132").Value2 = 20723.0316
$ws.Range("M132").Value2 = -18193.0316
$ws.Range("H136").Value2 = 941800.0600000001
$ws.Range("I136").Value2 = 35980.1
$ws.Range("J136").Value2 = 10000000
$ws.Range("K136").Value2 = 107940.3
$ws.Range("L136").Value2 = 30000000
$ws.Range("M136").Value2 = -105390.3
$ws.Range("N136").Value2 = -30005100
